$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 35600
$ws.Range("J87").Value = 34500
$ws.Range("L87").Value = 34500
$ws.Range("N87").Value = -36996
$ws.Range("H90").Value = 35600
$ws.Range("J90").Value = 34500
$ws.Range("L90").Value = 103500
$ws.Range("N90").Value = -115980
$ws.Range("H98").Value = 11877.421
$ws.Range("I98").Value = 11877.421
$ws.Range("K98").Value = 11877.421
$ws.Range("M98").Value = -10379.421
$ws.Range("H106").Value = 16702010
$ws.Range("I106").Value = 40280.383
$ws.Range("K106").Value = 40280.383
$ws.Range("M106").Value = -39649.383
$ws.Range("H122").Value = 11877.421
$ws.Range("I122").Value = 11877.421
$ws.Range("K122").Value = 35632.263
$ws.Range("M122").Value = -33182.263
$ws.Range("H129").Value = 54180.42
$ws.Range("J129").Value = 85472.664
$ws.Range("L129").Value = 256417.992
$ws.Range("N129").Value = -266417.992
$ws.Range("H132").Value = 3249192
$ws.Range("I132").Value = 3761759.2
$ws.Range("J132").Value = 2932.6667
$ws.Range("K132").Value = 11285277.6
$ws.Range("L132").Value = 8798.000100000001
$ws.Range("M132").Value = -11282747.6
$ws.Range("N132").Value = -13858.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1624.75
$ws.Range("I2").Value = 666.6667
$ws.Range("J2").Value = 2199.6
$ws.Range("K2").Value = 666.6667
$ws.Range("L2").Value = 2199.6
$ws.Range("M2").Value = -553.6667
$ws.Range("N2").Value = -2425.6
$ws.Range("H45").Value = 875.4286
$ws.Range("I45").Value = 819
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 819
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -442
$ws.Range("N45").Value = -1968
$ws.Range("H116").Value = 1624.75
$ws.Range("I116").Value = 666.6667
$ws.Range("J116").Value = 2199.6
$ws.Range("K116").Value = 666.6667
$ws.Range("L116").Value = 2199.6
$ws.Range("M116").Value = 1627.3333
$ws.Range("N116").Value = -6787.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1624.75
$ws.Range("I3").Value = 666.6667
$ws.Range("J3").Value = 2199.6
$ws.Range("K3").Value = 666.6667
$ws.Range("L3").Value = 2199.6
$ws.Range("M3").Value = -552.6667
$ws.Range("N3").Value = -2427.6
$ws.Range("H134").Value = 5745.943
$ws.Range("I134").Value = 7230.364
$ws.Range("J134").Value = 3233.8462
$ws.Range("K134").Value = 21691.092
$ws.Range("L134").Value = 9701.5386
$ws.Range("M134").Value = -19156.092
$ws.Range("N134").Value = -14771.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -830
$ws.Range("H24").Value = 1000
$ws.Range("I24").Value = 1000
$ws.Range("K24").Value = 1000
$ws.Range("M24").Value = -830
$ws.Range("H31").Value = 3127.6
$ws.Range("I31").Value = 2518
$ws.Range("J31").Value = 5130.5713
$ws.Range("K31").Value = 2518
$ws.Range("L31").Value = 5130.5713
$ws.Range("M31").Value = -2223
$ws.Range("N31").Value = -5720.5713
$ws.Range("H34").Value = 3127.6
$ws.Range("I34").Value = 2518
$ws.Range("J34").Value = 5130.5713
$ws.Range("K34").Value = 2518
$ws.Range("L34").Value = 5130.5713
$ws.Range("M34").Value = -2316
$ws.Range("N34").Value = -5534.5713
$ws.Range("H99").Value = 2218.9524
$ws.Range("I99").Value = 1936.1818
$ws.Range("J99").Value = 2530
$ws.Range("K99").Value = 1936.1818
$ws.Range("L99").Value = 2530
$ws.Range("M99").Value = -438.1818000000001
$ws.Range("N99").Value = -5526
$ws.Range("H126").Value = 2218.9524
$ws.Range("I126").Value = 1936.1818
$ws.Range("J126").Value = 2530
$ws.Range("K126").Value = 5808.5454
$ws.Range("L126").Value = 7590
$ws.Range("M126").Value = -3338.5454
$ws.Range("N126").Value = -12530
$ws.Range("H134").Value = 3556.2222
$ws.Range("I134").Value = 3629.257
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 10887.771
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -8352.771000000001
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1506.625
$ws.Range("I34").Value = 750
$ws.Range("J34").Value = 1960.6
$ws.Range("K34").Value = 2250
$ws.Range("L34").Value = 5881.799999999999
$ws.Range("M34").Value = -2166
$ws.Range("N34").Value = -6049.799999999999
$ws.Range("H39").Value = 3249.6875
$ws.Range("J39").Value = 3433
$ws.Range("L39").Value = 10299
$ws.Range("N39").Value = -10887
$ws.Range("H55").Value = 2692
$ws.Range("J55").Value = 2692
$ws.Range("L55").Value = 8076
$ws.Range("N55").Value = -8430
$ws.Range("H121").Value = 11545.789
$ws.Range("I121").Value = 12938
$ws.Range("J121").Value = 10533.272
$ws.Range("K121").Value = 38814
$ws.Range("L121").Value = 31599.816
$ws.Range("M121").Value = -37504
$ws.Range("N121").Value = -34219.81600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10332
$ws.Range("H113").Value = 17858394
$ws.Range("I113").Value = 50000904
$ws.Range("J113").Value = 1444.4445
$ws.Range("K113").Value = 50000904
$ws.Range("L113").Value = 1444.4445
$ws.Range("M113").Value = -49998734
$ws.Range("N113").Value = -5784.4445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71430104
$ws.Range("I7").Value = 1791.6666
$ws.Range("J7").Value = 500000000
$ws.Range("K7").Value = 1791.6666
$ws.Range("L7").Value = 500000000
$ws.Range("M7").Value = -1679.6666
$ws.Range("N7").Value = -500000224
$ws.Range("H40").Value = 2788.375
$ws.Range("I40").Value = 2758.1428
$ws.Range("K40").Value = 2758.1428
$ws.Range("M40").Value = -2622.1428
$ws.Range("H61").Value = 20835010
$ws.Range("I61").Value = 1555.4286
$ws.Range("J61").Value = 37038810
$ws.Range("K61").Value = 1555.4286
$ws.Range("L61").Value = 37038810
$ws.Range("M61").Value = -1353.4286
$ws.Range("N61").Value = -37039214
$ws.Range("H75").Value = 9578.5
$ws.Range("I75").Value = 9578.5
$ws.Range("K75").Value = 9578.5
$ws.Range("M75").Value = -8642.5
$ws.Range("H76").Value = 1500
$ws.Range("I76").Value = 1500
$ws.Range("K76").Value = 1500
$ws.Range("M76").Value = -1162
$ws.Range("H78").Value = 9578.5
$ws.Range("I78").Value = 9578.5
$ws.Range("K78").Value = 28735.5
$ws.Range("M78").Value = -24055.5
$ws.Range("H79").Value = 1500
$ws.Range("I79").Value = 1500
$ws.Range("K79").Value = 1500
$ws.Range("M79").Value = -330
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H113").Value = 20835010
$ws.Range("I113").Value = 1555.4286
$ws.Range("J113").Value = 37038810
$ws.Range("K113").Value = 1555.4286
$ws.Range("L113").Value = 37038810
$ws.Range("M113").Value = 614.5714
$ws.Range("N113").Value = -37043150
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 71430104
$ws.Range("I126").Value = 1791.6666
$ws.Range("J126").Value = 500000000
$ws.Range("K126").Value = 5374.9998
$ws.Range("L126").Value = 1500000000
$ws.Range("M126").Value = -2904.9998
$ws.Range("N126").Value = -1500004940
$ws.Range("H132").Value = 6202.85
$ws.Range("I132").Value = 7696.857
$ws.Range("J132").Value = 2716.8333
$ws.Range("K132").Value = 23090.571
$ws.Range("L132").Value = 8150.499899999999
$ws.Range("M132").Value = -20560.571
$ws.Range("N132").Value = -13210.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2226
$ws.Range("H5").Value = 37600000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H6").Value = 750
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -385
$ws.Range("N6").Value = -1230
